$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("CRP")
$ws4 = $wb.Worksheets.Item("CUL")
$ws5 = $wb.Worksheets.Item("LTW")
$ws6 = $wb.Worksheets.Item("WVR")

# --- ALC ---
$ws1.Range("H33").Value = 553.5714
$ws1.Range("I33").Value = 519.2308
$ws1.Range("J33").Value = 1000
$ws1.Range("K33").Value = 519.2308
$ws1.Range("L33").Value = 1000
$ws1.Range("M33").Value = -290.2308
$ws1.Range("N33").Value = -1458
$ws1.Range("H40").Value = 1084.5
$ws1.Range("I40").Value = 1023.82355
$ws1.Range("J40").Value = 1199.1111
$ws1.Range("K40").Value = 1023.82355
$ws1.Range("L40").Value = 1199.1111
$ws1.Range("M40").Value = -848.82355
$ws1.Range("N40").Value = -1549.1111
$ws1.Range("H51").Value = 3007.2144
$ws1.Range("I51").Value = 5000
$ws1.Range("J51").Value = 2853.923
$ws1.Range("K51").Value = 5000
$ws1.Range("L51").Value = 2853.923
$ws1.Range("M51").Value = -4516
$ws1.Range("N51").Value = -3821.923
$ws1.Range("H53").Value = 604.6
$ws1.Range("I53").Value = 264.91666
$ws1.Range("J53").Value = 918.1539
$ws1.Range("K53").Value = 264.91666
$ws1.Range("L53").Value = 918.1539
$ws1.Range("M53").Value = 372.08334
$ws1.Range("N53").Value = -2192.1539
$ws1.Range("H86").Value = 3175.7073
$ws1.Range("I86").Value = 2338.0688
$ws1.Range("J86").Value = 5200
$ws1.Range("K86").Value = 2338.0688
$ws1.Range("L86").Value = 5200
$ws1.Range("M86").Value = -1215.0688
$ws1.Range("N86").Value = -7446
$ws1.Range("H89").Value = 3175.7073
$ws1.Range("I89").Value = 2338.0688
$ws1.Range("J89").Value = 5200
$ws1.Range("K89").Value = 11690.344
$ws1.Range("L89").Value = 26000
$ws1.Range("M89").Value = -6074.344000000001
$ws1.Range("N89").Value = -37232
$ws1.Range("H98").Value = 2168.6135
$ws1.Range("I98").Value = 2525.4722
$ws1.Range("J98").Value = 562.75
$ws1.Range("K98").Value = 2525.4722
$ws1.Range("L98").Value = 562.75
$ws1.Range("M98").Value = -1027.4722
$ws1.Range("N98").Value = -3558.75
$ws1.Range("H106").Value = 1962.5
$ws1.Range("I106").Value = 1962.5
$ws1.Range("J106").Value = 0
$ws1.Range("K106").Value = 1962.5
$ws1.Range("L106").Value = 0
$ws1.Range("M106").Value = -1331.5
$ws1.Range("N106").Value = ""
$ws1.Range("H122").Value = 2168.6135
$ws1.Range("I122").Value = 2525.4722
$ws1.Range("J122").Value = 562.75
$ws1.Range("K122").Value = 7576.4166
$ws1.Range("L122").Value = 1688.25
$ws1.Range("M122").Value = -5126.4166
$ws1.Range("N122").Value = -6588.25
$ws1.Range("H137").Value = 46668710
$ws1.Range("I137").Value = 9805696
$ws1.Range("K137").Value = 29417088
$ws1.Range("M137").Value = -29414538
$ws1.Range("H138").Value = 1472871.8
$ws1.Range("I138").Value = 2223698.2
$ws1.Range("J138").Value = 3863.5652
$ws1.Range("K138").Value = 6671094.600000001
$ws1.Range("L138").Value = 11590.6956
$ws1.Range("M138").Value = -6665954.600000001
$ws1.Range("N138").Value = -21870.6956
$ws1.Range("H141").Value = 1066.2941
$ws1.Range("I141").Value = 766
$ws1.Range("J141").Value = 1787
$ws1.Range("K141").Value = 2298
$ws1.Range("L141").Value = 5361
$ws1.Range("M141").Value = 2882
$ws1.Range("N141").Value = -15721

# --- ARM ---
$ws2.Range("H74").Value = 6255163.5
$ws2.Range("I74").Value = 11765496
$ws2.Range("J74").Value = 10120.866
$ws2.Range("K74").Value = 11765496
$ws2.Range("L74").Value = 10120.866
$ws2.Range("M74").Value = -11764622
$ws2.Range("N74").Value = -11868.866
$ws2.Range("H77").Value = 6255163.5
$ws2.Range("I77").Value = 11765496
$ws2.Range("J77").Value = 10120.866
$ws2.Range("K77").Value = 58827480
$ws2.Range("L77").Value = 50604.33
$ws2.Range("M77").Value = -58823112
$ws2.Range("N77").Value = -59340.33
$ws2.Range("H132").Value = 1441.0294
$ws2.Range("I132").Value = 1415.3273
$ws2.Range("J132").Value = 1549.7693
$ws2.Range("K132").Value = 4245.9819
$ws2.Range("L132").Value = 4649.3079
$ws2.Range("M132").Value = -1715.9819
$ws2.Range("N132").Value = -9709.3079

# --- CRP ---
$ws3.Range("H7").Value = 106.77778
$ws3.Range("I7").Value = 90
$ws3.Range("J7").Value = 120.2
$ws3.Range("K7").Value = 90
$ws3.Range("L7").Value = 120.2
$ws3.Range("M7").Value = 23
$ws3.Range("N7").Value = -346.2
$ws3.Range("H22").Value = 268.54166
$ws3.Range("I22").Value = 180.61111
$ws3.Range("J22").Value = 532.3333
$ws3.Range("K22").Value = 180.61111
$ws3.Range("L22").Value = 532.3333
$ws3.Range("M22").Value = 169.38889
$ws3.Range("N22").Value = -1232.3333
$ws3.Range("H31").Value = 1517.6129
$ws3.Range("I31").Value = 1412.1
$ws3.Range("J31").Value = 1709.4546
$ws3.Range("K31").Value = 1412.1
$ws3.Range("L31").Value = 1709.4546
$ws3.Range("M31").Value = -1117.1
$ws3.Range("N31").Value = -2299.4546
$ws3.Range("H34").Value = 1517.6129
$ws3.Range("I34").Value = 1412.1
$ws3.Range("J34").Value = 1709.4546
$ws3.Range("K34").Value = 1412.1
$ws3.Range("L34").Value = 1709.4546
$ws3.Range("M34").Value = -1210.1
$ws3.Range("N34").Value = -2113.4546

# --- CUL ---
$ws4.Range("H5").Value = 970
$ws4.Range("I5").Value = 887.5
$ws4.Range("K5").Value = 2662.5
$ws4.Range("M5").Value = -2550.5
$ws4.Range("H57").Value = 865
$ws4.Range("I57").Value = 865
$ws4.Range("J57").Value = 0
$ws4.Range("K57").Value = 2595
$ws4.Range("L57").Value = 0
$ws4.Range("M57").Value = -2036
$ws4.Range("N57").Value = ""
$ws4.Range("H100").Value = 3891
$ws4.Range("I100").Value = 1980
$ws4.Range("J100").Value = 4528
$ws4.Range("K100").Value = 5940
$ws4.Range("L100").Value = 13584
$ws4.Range("M100").Value = -5129
$ws4.Range("N100").Value = -15206
$ws4.Range("H129").Value = 1129.6818
$ws4.Range("I129").Value = 699.6667
$ws4.Range("J129").Value = 1197.579
$ws4.Range("K129").Value = 2099.0001
$ws4.Range("L129").Value = 3592.737
$ws4.Range("M129").Value = 2900.9999
$ws4.Range("N129").Value = -13592.737
$ws4.Range("H131").Value = 4096.033
$ws4.Range("I131").Value = 4160.077
$ws4.Range("K131").Value = 12480.231
$ws4.Range("M131").Value = -7440.231
$ws4.Range("H135").Value = 970
$ws4.Range("I135").Value = 887.5
$ws4.Range("K135").Value = 7987.5
$ws4.Range("M135").Value = -5452.5

# --- LTW ---
$ws5.Range("H136").Value = 2557.8845
$ws5.Range("I136").Value = 1263.4375
$ws5.Range("J136").Value = 4629
$ws5.Range("K136").Value = 3790.3125
$ws5.Range("L136").Value = 13887
$ws5.Range("M136").Value = -1240.3125
$ws5.Range("N136").Value = -18987

# --- WVR ---
$ws6.Range("H136").Value = 9507.791999999999
$ws6.Range("I136").Value = 12739.941
$ws6.Range("K136").Value = 38219.823
$ws6.Range("M136").Value = -35669.823

Write-Output "Applied Ifrit_Profits edits"